$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# are forced to remain text (matching the source inlineStr cells) by temporarily
# switching the cell to a text number format, then restoring the default style.

$ws.Range('D2').Value = '60.722.75'
$ws.Range('E2').Value = '  +3.06%  '
$ws.Range('D3').Value = '2.690.12'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('D9').Value = '2.711.52'
$ws.Range('E9').Value = '  +2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.85%  '
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.341'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('E13').Value = '  +2.82%  '
$ws.Range('D14').Value = '3.163.39'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '60.702.77'
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.995.26'
$ws.Range('E16').Value = '  +13.20%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '21.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '347.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('E22').Value = '  +2.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.60'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.97%  '
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E26').Value = '  +6.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.993'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.55%  '
$ws.Range('D29').Value = '0.0₃0818'
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.07%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('E35').Value = '  +6.66%  '
$ws.Range('E36').Value = '  +9.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.941'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.882'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.35%  '
$ws.Range('E39').Value = '  +8.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '287.76'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0993'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.612'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('E45').Value = '  +2.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.995'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '2.136.52'
$ws.Range('E47').Value = '  +7.92%  '
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.26%  '
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('E51').Value = '  +1.75%  '
